$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 ("Camila") entirely, shifting row 3 ("Gonzalo") up to row 2
$ws.Rows.Item(2).Delete()
